$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (ALC)
$ws.Range("H11").Value = 41666900
$ws.Range("I11").Value = 41666900
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 41666900
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -41666760

# Row 100 (ALC)
$ws.Range("H100").Value = 127826.375
$ws.Range("I100").Value = 505002.5
$ws.Range("J100").Value = 2101
$ws.Range("K100").Value = 505002.5
$ws.Range("L100").Value = 2101
$ws.Range("M100").Value = -504461.5
$ws.Range("N100").Value = -3183

# Row 111 (ALC)
$ws.Range("H111").Value = 1501.5
$ws.Range("I111").Value = 1819.875
$ws.Range("J111").Value = 1319.5714
$ws.Range("K111").Value = 5459.625
$ws.Range("L111").Value = 3958.7142
$ws.Range("M111").Value = -2392.625
$ws.Range("N111").Value = -10092.7142

# Row 137 (ALC)
$ws.Range("H137").Value = 1495.225
$ws.Range("I137").Value = 1184.15
$ws.Range("J137").Value = 1806.3
$ws.Range("K137").Value = 3552.45
$ws.Range("L137").Value = 5418.9
$ws.Range("M137").Value = -1002.45
$ws.Range("N137").Value = -10518.9

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM)
$ws.Range("H94").Value = 553.73334
$ws.Range("I94").Value = 518
$ws.Range("J94").Value = 696.6667
$ws.Range("K94").Value = 518
$ws.Range("L94").Value = 696.6667
$ws.Range("M94").Value = -67
$ws.Range("N94").Value = -1598.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 9265.302
$ws.Range("I31").Value = 3290.372
$ws.Range("J31").Value = 22111.4
$ws.Range("K31").Value = 3290.372
$ws.Range("L31").Value = 22111.4
$ws.Range("M31").Value = -2995.372
$ws.Range("N31").Value = -22701.4

# Row 34 (CRP)
$ws.Range("H34").Value = 9265.302
$ws.Range("I34").Value = 3290.372
$ws.Range("J34").Value = 22111.4
$ws.Range("K34").Value = 3290.372
$ws.Range("L34").Value = 22111.4
$ws.Range("M34").Value = -3088.372
$ws.Range("N34").Value = -22515.4

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 2714.8484
$ws.Range("I80").Value = 2250
$ws.Range("J80").Value = 2744.8386
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 2744.8386
$ws.Range("M80").Value = -1252
$ws.Range("N80").Value = -4740.8386

# Row 83 (GSM)
$ws.Range("H83").Value = 2714.8484
$ws.Range("I83").Value = 2250
$ws.Range("J83").Value = 2744.8386
$ws.Range("K83").Value = 11250
$ws.Range("L83").Value = 13724.193
$ws.Range("M83").Value = -6258
$ws.Range("N83").Value = -23708.193

# Row 102 (GSM)
$ws.Range("H102").Value = 1093.4
$ws.Range("I102").Value = 1016.75
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 1016.75
$ws.Range("L102").Value = 1400
$ws.Range("M102").Value = 605.25
$ws.Range("N102").Value = -4644

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2470.111
$ws.Range("I7").Value = 2266.4
$ws.Range("J7").Value = 2724.75
$ws.Range("K7").Value = 2266.4
$ws.Range("L7").Value = 2724.75
$ws.Range("M7").Value = -2154.4
$ws.Range("N7").Value = -2948.75

# Row 16 (LTW)
$ws.Range("H16").Value = 12812.375
$ws.Range("I16").Value = 357
$ws.Range("J16").Value = 100000
$ws.Range("K16").Value = 357
$ws.Range("L16").Value = 100000
$ws.Range("M16").Value = -187
$ws.Range("N16").Value = -100340

# Row 40 (LTW)
$ws.Range("H40").Value = 1180
$ws.Range("I40").Value = 1180
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1180
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1044

# Row 61 (LTW)
$ws.Range("H61").Value = 3970
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 3970
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 3970
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -4374

# Row 87 (LTW)
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()

# Row 90 (LTW)
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()

# Row 110 (LTW)
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 113 (LTW)
$ws.Range("H113").Value = 3970
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3970
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3970
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8310

# Row 122 (LTW)
$ws.Range("H122").Value = 62505348
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 62505348
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 187516044
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -187520944

# Row 123 (LTW)
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 126 (LTW)
$ws.Range("H126").Value = 2470.111
$ws.Range("I126").Value = 2266.4
$ws.Range("J126").Value = 2724.75
$ws.Range("K126").Value = 6799.200000000001
$ws.Range("L126").Value = 8174.25
$ws.Range("M126").Value = -4329.200000000001
$ws.Range("N126").Value = -13114.25

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 2329.7144
$ws.Range("I81").Value = 547.5
$ws.Range("J81").Value = 3042.6
$ws.Range("K81").Value = 1095
$ws.Range("L81").Value = 6085.2
$ws.Range("M81").Value = -34
$ws.Range("N81").Value = -8207.200000000001

# Row 84 (WVR)
$ws.Range("H84").Value = 2329.7144
$ws.Range("I84").Value = 547.5
$ws.Range("J84").Value = 3042.6
$ws.Range("K84").Value = 5475
$ws.Range("L84").Value = 30426
$ws.Range("M84").Value = -171
$ws.Range("N84").Value = -41034

# Row 96 (WVR)
$ws.Range("H96").Value = 3000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 3000
$ws.Range("N96").Value = -5746

# Row 107 (WVR)
$ws.Range("H107").Value = 13513927
$ws.Range("I107").Value = 408.85715
$ws.Range("J107").Value = 21739548
$ws.Range("K107").Value = 1226.57145
$ws.Range("L107").Value = 65218644
$ws.Range("M107").Value = 693.4285500000001
$ws.Range("N107").Value = -65222484

# Row 110 (WVR)
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 111 (WVR)
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 114 (WVR)
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 115 (WVR)
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# Row 122 (WVR)
$ws.Range("H122").Value = 4235.276
$ws.Range("I122").Value = 4448.2593
$ws.Range("J122").Value = 1360
$ws.Range("K122").Value = 13344.7779
$ws.Range("L122").Value = 4080
$ws.Range("M122").Value = -10894.7779

# Row 123 (WVR)
$ws.Range("H123").Value = 20429
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20429
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20429
$ws.Range("N123").Value = -30229

# Row 125 (WVR)
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 126 (WVR)
$ws.Range("H126").Value = 111111990
$ws.Range("I126").Value = 142857780
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 428573340
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -428570870
$ws.Range("N126").Value = -10190
